$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header text changes (case -> lowercase) ---
$ws.Range("A1").Value = "names"
$ws.Range("B1").Value = "careers"
$ws.Range("C1").Value = "base"

# --- New column D header/content ---
$ws.Range("D2").Value = "job-title"

# --- New row 3 (Jump Trading) ---
$ws.Range("A3").Value = "Jump Trading"
$ws.Range("B3").Value = "https://www.jumptrading.com/careers/"
$ws.Range("C3").Value = "https://www.jumptrading.com/"

# --- New column E header/content ---
$ws.Range("E2").Value = "job-location-name"

# --- Row 3 class strings (D3/E3) ---
$ws.Range("D3").Value = "text-xl lg:text-2xl font-medium text-black group-hover:text-jump-red"
$ws.Range("E3").Value = "text-base lg:text-lg text-dark-gray group-hover:text-black"

# --- Row 1 new headers D1/E1 ---
$ws.Range("D1").Value = "job_classes"
$ws.Range("E1").Value = "location_classes"

# --- Formatting: wrap text on D3/E3 (forces a dedicated font + style, matching the
#     new cellXfs entry created when the original author applied this formatting) ---
$ws.Range("D3:E3").Font.Name = "Calibri"
$ws.Range("D3:E3").WrapText = $true

# --- Row height for row 3 ---
$ws.Rows.Item(3).RowHeight = 59

# --- Column widths for D and E ---
$ws.Columns.Item(4).ColumnWidth = 20.5
$ws.Columns.Item(5).ColumnWidth = 18.1666667

# --- Selection to match target ---
$ws.Range("D11").Select()
